# edit.ps1 - apply the tracked changes to twoRobotRegionV.pptx
#
# Summary of changes (per the OOXML diff):
#  1. Every cached "datetimeFigureOut" footer field ("4/17/18") on the slide
#     master and all 11 slide layouts is refreshed to "5/22/18".
#  2. The translucent highlight rectangle ("Rectangle 9") on slide 1 switches
#     its solid fill from a theme-tinted accent3 color to a flat green
#     (RGB 00FF00) with ~22.7% opacity.
#  3. The rotated "TextBox 12" caption on slide 1 grows taller (its rotated
#     bounding box extends further) and its wording changes from
#     "Reachable set for s_1" to "2-move reachable set for s_1" at a smaller
#     font size (48pt -> 40pt), keeping the "s" italic and the trailing "1"
#     as a subscript.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached date field text on the master + every layout
# ---------------------------------------------------------------------
function Update-DateText($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "4/17/18") {
                $sh.TextFrame.TextRange.Text = "5/22/18"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateText $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateText $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 1: recolor the highlight rectangle's fill
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$rect = $slide.Shapes.Item(3)   # "Rectangle 9"

$rect.Fill.ForeColor.RGB = 65280      # RGB(0,255,0) -> 00FF00
$rect.Fill.Transparency = 0.77255     # alpha 22745/100000 ~= 22.745%

# ---------------------------------------------------------------------
# 3) Slide 1: resize + retext the rotated "TextBox 12" caption
# ---------------------------------------------------------------------
$box = $slide.Shapes.Item(6)   # "TextBox 12"

# Grow the (rotated) bounding box: Left/Width stay put, Top moves up and
# Height grows so the box still reads correctly once rotated 270 degrees.
$box.Top = 328.07708661417325
$box.Height = 152.67653743307085

$tr = $box.TextFrame.TextRange

# Shrink all the existing runs from 48pt to 40pt first (this edits the
# pre-existing runs in place, preserving their dirty/formatting markers).
$tr.Font.Size = 40

# Retype "Reachable " (the word plus its trailing space) as
# "2-move reachable " so the sentence now reads
# "2-move reachable set for s1".
$prefix = $tr.Characters(1, 10)
$prefix.Text = "2-move reachable "
